$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- Step 1: insert a new column before G, duplicating column F ----
$ws1.Columns.Item(6).Copy()
$ws1.Columns.Item(7).Insert()
$ws1.Columns.Item(7).ColumnWidth = $ws1.Columns.Item(6).ColumnWidth()

# ---- Step 2: append the five new test-data rows (17-21) on Sheet1 ----
$newRows = @(
    @("TC_Register_003", "female", "Anusha",     "MK",        "anumk@gmail.com",     "anumk#9669"),
    @("TC_Register_003", "male",   "renuka",     "prasad",    "renuka@gmail.com",    "prasad#123"),
    @("TC_Register_003", "male",   "pradyumna ", "r",         "pradyumna@gmail.com", "r#123"),
    @("TC_Register_003", "male",   "vikas",      "r",         "vikas@gmail.com",     "vikas#123"),
    @("TC_Register_003", "male",   "prem",       "choudhary", "prem@gmail.com",      "choudhary#123")
)

$r = 17
foreach ($row in $newRows) {
    $ws1.Range("A$r").Value = $row[0]
    $ws1.Range("B$r").Value = $row[1]
    $ws1.Range("C$r").Value = $row[2]
    $ws1.Range("D$r").Value = $row[3]
    $ws1.Range("E$r").Value = $row[4]
    $ws1.Hyperlinks.Add($ws1.Range("E$r"), "mailto:" + $row[4])
    $ws1.Range("F$r").Value = $row[5]
    $ws1.Range("G$r").Value = $row[5]
    $r = $r + 1
}

# restore the shared "Hyperlink" cell style on the E column of the new rows
# (Hyperlinks.Add applies a private per-call style variant otherwise)
$ws1.Range("E17:E21").Style = "Hyperlink"

# ---- Step 3: add Sheet2 (after Sheet1) with the TC_Register_003/Anusha row ----
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "TC_Register_003"
$ws2.Range("B1").Value = "female"
$ws2.Range("C1").Value = "Anusha"
$ws2.Range("D1").Value = "MK"
$ws2.Range("E1").Value = "anumk@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("E1"), "mailto:anumk@gmail.com")
$ws2.Range("F1").Value = "anumk#9669"
$ws2.Range("G1").Value = "anumk#9669"
$ws2.Range("E1").Style = "Hyperlink"

$ws2.Columns.Item(5).ColumnWidth = 16.94
$ws2.Columns.Item(6).ColumnWidth = 14.0
$ws2.Columns.Item(7).ColumnWidth = 12.7

# ---- Step 4: restore cursor/selection state to match the saved workbook ----
$null = $ws2.Range("G14").Select()
$null = $ws1.Range("F24").Select()
$ws1.Activate()

Write-Host "done"
